# Updates cryptos list prices / volume(1h) percentages (commit: Updated cryptos list)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.988.42"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.871.04"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.37%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5079"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3936"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08185"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.41%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.094"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.76"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("D13").Value = "1.863.63"
$ws.Range("E13").Value = "  -3.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.270"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.173"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.92"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.16%  "
$ws.Range("E18").Value = "  -3.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06426"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.69%  "
$ws.Range("E20").Value = "  -1.86%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "29.977.36"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.805"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.22%  "
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.145"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.50%  "
$ws.Range("D26").Value = "2.082.72"
$ws.Range("E26").Value = "  -3.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.10"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("E28").Value = "  -1.44%  "
$ws.Range("E29").Value = "  -9.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.05"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.057"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.901"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.741"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02418"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.250"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06338"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2142"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.46%  "
$ws.Range("E39").Value = "  -5.48%  "
$ws.Range("E40").Value = "  -5.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6302"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.22"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.196"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.31%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5902"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.91"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.88%  "
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.995"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.67"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("E50").Value = "  -3.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.128"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.07%  "
